$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates ---------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 1128860
# Cant. Trabajadores (worker count) now that one worker's rows are removed
$ws.Range("C13").Value = 1

# --- Re-sort the "Periodo Mora" detail rows (16-56) ascending ---------
# Previously the 41 rows for OTTO MARIO OSPINA BERRIO (periods 1705..2009)
# were listed most-recent-first (2009 down to 1705). They are now listed
# oldest-first (1705 up to 2009), and the "Valor Mora" amount that used to
# sit on the first row now sits on the last row.
$periods = @("1705","1706","1707","1708","1709","1710","1711","1712", `
             "1801","1802","1803","1804","1805","1806","1807","1808", `
             "1809","1810","1811","1812","1901","1902","1903","1904", `
             "1905","1906","1907","1908","1909","1910","1911","1912", `
             "2001","2002","2003","2004","2005","2006","2007","2008","2009")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 5).Value = $periods[$i]
}
$ws.Cells.Item(16, 6).Value = 27578
$ws.Cells.Item(56, 6).Value = 25740

# --- Remove the second worker (JAIME ALFONSO HOLLMAN GONZALEZ) --------
# His three mora rows are deleted outright; remaining rows below shift up.
$ws.Rows("57:59").Delete()

# Give the new last detail row (56) the heavier bottom border that used
# to close out the table (it previously sat on row 59).
$ws.Range("B56:J56").Borders.Item(9).LineStyle = 1
$ws.Range("B56:J56").Borders.Item(9).Weight = 2
$ws.Range("B56:J56").Borders.Item(9).ColorIndex = 1
